# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for column G ("K") for data rows 2-42
$newK = @{
    2  = 4
    3  = 2
    4  = 3
    5  = 8
    6  = 8
    7  = 7
    8  = 0
    9  = 8
    10 = 9
    11 = 8
    12 = 9
    13 = 8
    14 = 7
    15 = 6
    16 = 10
    17 = 4
    18 = 7
    19 = 8
    20 = 7
    21 = 9
    22 = 6
    23 = 6
    24 = 4
    25 = 4
    26 = 5
    27 = 7
    28 = 7
    29 = 4
    30 = 6
    31 = 4
    32 = 5
    33 = 6
    34 = 3
    35 = 11
    36 = 3
    37 = 8
    38 = 4
    39 = 3
    40 = 5
    41 = 3
    42 = 2
}

foreach ($row in $newK.Keys) {
    $ws.Cells.Item($row, 7).Value = $newK[$row]
}
